$d = $word.ActiveDocument

# The document currently starts with two paragraphs:
#   1) "Poverty is the Face of Christ"   (pStyle Heading1, single run)
#   2) "By Dorothy Day"                  (no pStyle, single bold run)
#
# The target (pandoc-style title block) turns these into:
#   1) pStyle "Title" paragraph with the title split word-by-word into
#      separate runs (with separate single-space runs between words).
#   2) pStyle "Authors" paragraph ("By " prefix and bold removed) with the
#      author's name likewise split word-by-word into separate runs.
#
# Strategy: insert the two brand-new, fully-formed paragraphs at the very
# start of the document (a collapsed range at position 0), which cleanly
# pushes the existing two paragraphs down without disturbing them. Then
# delete the (now-shifted) original two paragraphs outright. Building the
# replacement paragraphs from scratch this way (rather than overwriting the
# old ranges in place) avoids picking up any stray formatting/ids from the
# paragraphs being replaced.

$newParagraphsXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Title"/></w:pPr>
<w:r><w:t xml:space="preserve">Poverty</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">is</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">the</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Face</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">of</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Christ</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:pStyle w:val="Authors"/></w:pPr>
<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">Day</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($newParagraphsXml)

# The two original paragraphs have now been pushed down to indices 3 and 4
# (1-based): re-fetch them fresh and delete the pair (including their
# paragraph marks) in one go.
$oldTitlePara = $d.Paragraphs.Item(3)
$oldByLinePara = $d.Paragraphs.Item(4)
$oldRange = $d.Range($oldTitlePara.Range.Start, $oldByLinePara.Range.End)
$oldRange.Delete()
